$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Collect paragraphs to remove completely:
#   1) Paragraphs that contain an inline image (the 3 screenshots)
#   2) Paragraphs whose only content is the long dash separator line
#   3) Empty body paragraphs (not inside a table) used purely as
#      spacer paragraphs before a following block (pPr/spacing
#      w:before="40" -> SpaceBefore = 2pt)
# -----------------------------------------------------------------

$targets = @()

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $clean = $text.TrimEnd([char]13)

    $hasShape = $p.Range.InlineShapes.Count -gt 0
    $isSeparator = $clean.StartsWith("────")
    $inTable = $p.Range.Information(12)
    $isEmptySpacer = ($clean.Length -eq 0) -and (-not $hasShape) -and (-not $inTable) -and ($p.Format.SpaceBefore -eq 2)

    if ($hasShape -or $isSeparator -or $isEmptySpacer) {
        $targets += $p
    }
}

# Delete from the end of the document backwards so earlier ranges
# stay valid while we work.
for ($i = $targets.Count - 1; $i -ge 0; $i--) {
    $targets[$i].Range.Delete()
}
